# Applies two kinds of changes to the Chelsea stats workbook:
#  1. Renames the stat sheets to have human-friendly, spaced-out names
#     (and "Goal & Shot Creation" / "Miscellaneous Stats" full names).
#  2. Bumps the "Age" column (column E, formatted as "YY-DDD") of every
#     player row on every stats sheet forward by one day, leaving the
#     "Squad Total" / "Opponent Total" summary rows (which hold an average
#     like "24.3") untouched.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------
$renameMap = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"    = "Shooting Stats"
    "PassingStats"     = "Passing Stats"
    "PassTypes"        = "Pass Types"
    "GoalShotCreation" = "Goal & Shot Creation"
    "DefensiveActions" = "Defensive Actions"
    "PlayingTime"      = "Playing Time"
    "MiscStats"        = "Miscellaneous Stats"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renameMap.ContainsKey($oldName)) {
        $ws.Name = $renameMap[$oldName]
    }
}

# --- 2. Increment the day portion of the Age column on every stats sheet
foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Matches") {
        continue
    }

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 4; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        $val = $cell.Value2

        if ($val -match '^(\d{2})-(\d{3})$') {
            $years = $matches[1]
            $day = [int]$matches[2] + 1
            $dayStr = "$day".PadLeft(3, '0')
            $cell.Value2 = "$years-$dayStr"
        }
    }
}
